# Generate Report for Handoff
#
# File "b" (b.md) has gone through a new handoff cycle: its status moves
# from "Handed back: in sync with en-US" to "Ready for handoff", a fresh
# handoff file/datetime is recorded for each locale, and the Overview
# sheet / per-locale sheets are updated to reflect it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the b.md row.
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-28-12 20:28:17"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 (b.md) gets a new handoff status/file/datetime.
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-12 20:28:14"

# The hyperlink collection only supports append, not in-place edit, so
# rebuild it from scratch with the updated display text for D3 (the
# underlying rIds / target URLs are unchanged, matching the original
# workbook's stale-link behaviour).
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/02f0ed4076698f83c95420b8b7fbf277ad050abc/e2e/a.md", "", "", "a.md")
$zh.Hyperlinks.Add($zh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/02f0ed4076698f83c95420b8b7fbf277ad050abc/e2e/a.md", "", "", ".md")
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f7860dc39fba0f5efdb4a2e8647298cbddeb41e1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/2de1bf60cedf9ab7d2c6682f911f40bd213fba08/e2e/a.md", "", "", "a.md")
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c2171a740335389d182ddac5f201504c2e5e81f5/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/02f0ed4076698f83c95420b8b7fbf277ad050abc/e2e/b.md", "", "", "b.md")
$zh.Hyperlinks.Add($zh.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/02f0ed4076698f83c95420b8b7fbf277ad050abc/e2e/b.md", "", "", ".md")
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f7860dc39fba0f5efdb4a2e8647298cbddeb41e1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/2de1bf60cedf9ab7d2c6682f911f40bd213fba08/e2e/a.md", "", "", "a.md")
$zh.Hyperlinks.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c2171a740335389d182ddac5f201504c2e5e81f5/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")

# ---------------------------------------------------------------------
# de-de sheet: row 3 (b.md) gets a new handoff status/file/datetime.
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$de.Range("E3").Value = "2016-03-12 20:28:17"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/02f0ed4076698f83c95420b8b7fbf277ad050abc/e2e/a.md", "", "", "a.md")
$de.Hyperlinks.Add($de.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/02f0ed4076698f83c95420b8b7fbf277ad050abc/e2e/a.md", "", "", ".md")
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/828cced5336fd47a2c7cfb6139635cf863f36f3d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/af9abaed3782b864fc8e45bc6c4d149aa94033eb/e2e/a.md", "", "", "a.md")
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/de2c3c674dc27f05eeb34a543be2f0558b7cade2/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/02f0ed4076698f83c95420b8b7fbf277ad050abc/e2e/b.md", "", "", "b.md")
$de.Hyperlinks.Add($de.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/02f0ed4076698f83c95420b8b7fbf277ad050abc/e2e/b.md", "", "", ".md")
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/828cced5336fd47a2c7cfb6139635cf863f36f3d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/af9abaed3782b864fc8e45bc6c4d149aa94033eb/e2e/a.md", "", "", "a.md")
$de.Hyperlinks.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/de2c3c674dc27f05eeb34a543be2f0558b7cade2/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
